$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "43.877.50"
$ws.Range("E2").Value = "  +0.45%  "
$ws.Range("D3").Value = "2.357.86"
$ws.Range("E3").Value = "  +0.40%  "
$ws.Range("E4").Value = "  +0.26%  "
$ws.Range("D5").Value = "'0.671"
$ws.Range("E5").Value = "  +3.28%  "
$ws.Range("D6").Value = "'235.85"
$ws.Range("E6").Value = "  +0.86%  "
$ws.Range("D7").Value = "'73.31"
$ws.Range("E7").Value = "  +10.84%  "
$ws.Range("E8").Value = "  +0.10%  "
$ws.Range("D9").Value = "'0.562"
$ws.Range("E9").Value = "  +23.75%  "
$ws.Range("D10").Value = "'0.0988"
$ws.Range("E10").Value = "  +1.42%  "
$ws.Range("D11").Value = "'28.16"
$ws.Range("E11").Value = "  +4.70%  "
$ws.Range("E12").Value = "  +1.90%  "
$ws.Range("D13").Value = "2.709.15"
$ws.Range("E13").Value = "  +0.45%  "
$ws.Range("D14").Value = "'16.73"
$ws.Range("E14").Value = "  +8.20%  "
$ws.Range("D15").Value = "'6.77"
$ws.Range("E15").Value = "  +9.55%  "
$ws.Range("D16").Value = "'0.888"
$ws.Range("E16").Value = "  +3.46%  "
$ws.Range("D17").Value = "2.378.99"
$ws.Range("E17").Value = "  +1.30%  "
$ws.Range("D18").Value = "43.839.56"
$ws.Range("E18").Value = "  +0.31%  "
$ws.Range("E19").Value = "  +3.01%  "
$ws.Range("D20").Value = "'78.18"
$ws.Range("E20").Value = "  +5.54%  "
$ws.Range("E21").Value = "  +2.19%  "
$ws.Range("D22").Value = "'253.90"
$ws.Range("E22").Value = "  +1.73%  "
$ws.Range("E23").Value = "  -0.02%  "
$ws.Range("E24").Value = "  -1.23%  "
$ws.Range("E25").Value = "  +2.61%  "
$ws.Range("D26").Value = "'10.67"
$ws.Range("E26").Value = "  +7.24%  "
$ws.Range("E27").Value = "  +0.57%  "
$ws.Range("D28").Value = "'22.48"
$ws.Range("E28").Value = "  +0.25%  "
$ws.Range("D29").Value = "'172.76"
$ws.Range("E29").Value = "  -1.22%  "
$ws.Range("E30").Value = "  +9.89%  "
$ws.Range("E31").Value = "  +0.48%  "
$ws.Range("E32").Value = "  +5.63%  "
$ws.Range("D33").Value = "'5.20"
$ws.Range("E33").Value = "  +3.97%  "
$ws.Range("E34").Value = "  +4.49%  "
$ws.Range("D35").Value = "'5.21"
$ws.Range("E35").Value = "  +3.69%  "
$ws.Range("E36").Value = "  +1.77%  "
$ws.Range("D37").Value = "'2.44"
$ws.Range("E37").Value = "  -0.01%  "
$ws.Range("D38").Value = "'6.41"
$ws.Range("E38").Value = "  -2.51%  "
$ws.Range("D39").Value = "'0.0270"
$ws.Range("E39").Value = "  +6.69%  "
$ws.Range("D40").Value = "'19.23"
$ws.Range("E40").Value = "  +7.04%  "
$ws.Range("D41").Value = "'8.96"
$ws.Range("E41").Value = "  -2.65%  "
$ws.Range("E42").Value = "  +0.18%  "
$ws.Range("B43").Value = "Algorand"
$ws.Range("C43").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D43").Value = "'0.186"
$ws.Range("E43").Value = "  +13.69%  "
$ws.Range("B44").Value = "ARBITRUM"
$ws.Range("C44").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D44").Value = "'1.16"
$ws.Range("E44").Value = "  -2.47%  "
$ws.Range("D45").Value = "'0.0974"
$ws.Range("E45").Value = "  +2.08%  "
$ws.Range("E46").Value = "  +1.85%  "
$ws.Range("D47").Value = "'4.44"
$ws.Range("E47").Value = "  +1.60%  "
$ws.Range("D48").Value = "'97.72"
$ws.Range("E48").Value = "  -1.95%  "
$ws.Range("B49").Value = "Maker"
$ws.Range("C49").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D49").Value = "1.437.45"
$ws.Range("E49").Value = "  -0.81%  "
$ws.Range("B50").Value = "NEARProtocol"
$ws.Range("C50").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D50").Value = "'2.31"
$ws.Range("E50").Value = "  +0.25%  "
$ws.Range("D51").Value = "2.582.14"
$ws.Range("E51").Value = "  +0.43%  "
